$d = $word.ActiveDocument

# 1) Insert a new paragraph BEFORE the paragraph that starts with
#    "* Dificultades con la identificación y operación con términos semejantes."
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("* Dificultades con la identificación y operación con términos semejantes.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $insertPoint1 = $d.Range($rng1.Start, $rng1.Start)
    $insertPoint1.InsertBefore("* No se conoce o no se utiliza el criterio de la balanza. `r")
}

# 2) Insert a new paragraph AFTER the paragraph that contains
#    "* Dificultades con identificar y operar con fracciones algebraicas."
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("* Dificultades con identificar y operar con fracciones algebraicas.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $insertPoint2 = $d.Range($rng2.End, $rng2.End)
    $insertPoint2.InsertAfter("`r* Variabilidad de las respuestas en un grupo, cuando el procedimiento es de 1 o 2 páginas. ")
}
